$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.138.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.896.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5200'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3765'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07284'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9003'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08215'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.942.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.324'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008621'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.172.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.094'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.424'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.310'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.744'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.806'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.857'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09208'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7960'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.217'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.430'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.958'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.607'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5721'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02000'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.078'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.012'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.557'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.28'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.96%  '
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4866'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.619'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.42'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05930'
$ws.Range("D51").Style = "Normal"
